$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 553
$ws.Range("I2").Value = 553
$ws.Range("K2").Value = 553
$ws.Range("M2").Value = -440
$ws.Range("H5").Value = 1213.625
$ws.Range("I5").Value = 161
$ws.Range("K5").Value = 161
$ws.Range("M5").Value = -46
$ws.Range("H9").Value = 206.53334
$ws.Range("I9").Value = 197.33333
$ws.Range("J9").Value = 243.33333
$ws.Range("K9").Value = 197.33333
$ws.Range("L9").Value = 243.33333
$ws.Range("M9").Value = -28.33332999999999
$ws.Range("N9").Value = -581.3333299999999
$ws.Range("H118").Value = 1671
$ws.Range("J118").Value = 2996.6667
$ws.Range("L118").Value = 8990.000100000001
$ws.Range("N118").Value = -12304.0001
$ws.Range("H127").Value = 4262.091
$ws.Range("I127").Value = 1937.6
$ws.Range("J127").Value = 6199.1665
$ws.Range("K127").Value = 5812.799999999999
$ws.Range("L127").Value = 18597.4995
$ws.Range("M127").Value = -852.7999999999993
$ws.Range("N127").Value = -28517.4995
$ws.Range("H138").Value = 4438.797
$ws.Range("I138").Value = 2532.8333
$ws.Range("J138").Value = 4878.635
$ws.Range("K138").Value = 7598.499899999999
$ws.Range("L138").Value = 14635.905
$ws.Range("M138").Value = -2458.499899999999
$ws.Range("N138").Value = -24915.905

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1504.1666
$ws.Range("I5").Value = 1459.091
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 1459.091
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = -1347.091
$ws.Range("N5").Value = -2224
$ws.Range("H32").Value = 6164.6514
$ws.Range("I32").Value = 6079.953
$ws.Range("K32").Value = 6079.953
$ws.Range("M32").Value = -5792.953
$ws.Range("H45").Value = 7638.3076
$ws.Range("I45").Value = 7716.5
$ws.Range("K45").Value = 7716.5
$ws.Range("M45").Value = -7339.5
$ws.Range("H46").Value = 10047
$ws.Range("J46").Value = 10047
$ws.Range("L46").Value = 10047
$ws.Range("N46").Value = -10685
$ws.Range("H63").Value = 5490.5713
$ws.Range("I63").Value = 3858.5
$ws.Range("J63").Value = 7666.6665
$ws.Range("K63").Value = 3858.5
$ws.Range("L63").Value = 7666.6665
$ws.Range("M63").Value = -3172.5
$ws.Range("N63").Value = -9038.666499999999
$ws.Range("H66").Value = 5490.5713
$ws.Range("I66").Value = 3858.5
$ws.Range("J66").Value = 7666.6665
$ws.Range("K66").Value = 19292.5
$ws.Range("L66").Value = 38333.3325
$ws.Range("M66").Value = -15860.5
$ws.Range("N66").Value = -45197.3325
$ws.Range("H92").Value = 220034700
$ws.Range("J92").Value = 275032000
$ws.Range("L92").Value = 275032000
$ws.Range("N92").Value = -275036992
$ws.Range("H94").Value = 220029170
$ws.Range("J94").Value = 264020000
$ws.Range("L94").Value = 264020000
$ws.Range("N94").Value = -264021802
$ws.Range("H132").Value = 6263.1587
$ws.Range("I132").Value = 6154.115
$ws.Range("J132").Value = 6778.636
$ws.Range("K132").Value = 18462.345
$ws.Range("L132").Value = 20335.908
$ws.Range("M132").Value = -15932.345
$ws.Range("N132").Value = -25395.908
$ws.Range("H140").Value = 74999.336
$ws.Range("J140").Value = 74999.336
$ws.Range("L140").Value = 74999.336
$ws.Range("N140").Value = -85359.336

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1504.1666
$ws.Range("I4").Value = 1459.091
$ws.Range("J4").Value = 2000
$ws.Range("K4").Value = 1459.091
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = -1344.091
$ws.Range("N4").Value = -2230
$ws.Range("H134").Value = 17322.777
$ws.Range("I134").Value = 20129.357
$ws.Range("K134").Value = 60388.071
$ws.Range("M134").Value = -57853.071

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2500
$ws.Range("I16").Value = 2833.3333
$ws.Range("K16").Value = 2833.3333
$ws.Range("M16").Value = -2546.3333
$ws.Range("H22").Value = 660.8889
$ws.Range("I22").Value = 362.25
$ws.Range("J22").Value = 899.8
$ws.Range("K22").Value = 362.25
$ws.Range("L22").Value = 899.8
$ws.Range("M22").Value = -12.25
$ws.Range("N22").Value = -1599.8
$ws.Range("H31").Value = 2717.541
$ws.Range("I31").Value = 2190.0576
$ws.Range("K31").Value = 2190.0576
$ws.Range("M31").Value = -1895.0576
$ws.Range("H34").Value = 2717.541
$ws.Range("I34").Value = 2190.0576
$ws.Range("K34").Value = 2190.0576
$ws.Range("M34").Value = -1988.0576
$ws.Range("H86").Value = 14000
$ws.Range("I86").Value = 10000
$ws.Range("J86").Value = 22000
$ws.Range("K86").Value = 10000
$ws.Range("L86").Value = 22000
$ws.Range("M86").Value = -8877
$ws.Range("N86").Value = -24246
$ws.Range("H89").Value = 14000
$ws.Range("I89").Value = 10000
$ws.Range("J89").Value = 22000
$ws.Range("K89").Value = 50000
$ws.Range("L89").Value = 110000
$ws.Range("M89").Value = -44384
$ws.Range("N89").Value = -121232
$ws.Range("H113").Value = 2500
$ws.Range("I113").Value = 2833.3333
$ws.Range("K113").Value = 2833.3333
$ws.Range("M113").Value = -663.3332999999998
$ws.Range("H116").Value = 69999
$ws.Range("J116").Value = 69999
$ws.Range("L116").Value = 69999
$ws.Range("N116").Value = -79177

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 13355821
$ws.Range("I4").Value = 690514.6
$ws.Range("J4").Value = 105179290
$ws.Range("K4").Value = 2071543.8
$ws.Range("L4").Value = 315537870
$ws.Range("M4").Value = -2071431.8
$ws.Range("N4").Value = -315538094
$ws.Range("H39").Value = 3085.68
$ws.Range("J39").Value = 5292.2856
$ws.Range("L39").Value = 15876.8568
$ws.Range("N39").Value = -16464.8568
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("N92").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 28995.643
$ws.Range("I126").Value = 45198.6
$ws.Range("K126").Value = 135595.8
$ws.Range("M126").Value = -133125.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 11665
$ws.Range("I20").Value = 11665
$ws.Range("K20").Value = 11665
$ws.Range("M20").Value = -11439
$ws.Range("H22").Value = 6407.706
$ws.Range("I22").Value = 7497.95
$ws.Range("J22").Value = 4850.2144
$ws.Range("K22").Value = 7497.95
$ws.Range("L22").Value = 4850.2144
$ws.Range("M22").Value = -7202.95
$ws.Range("N22").Value = -5440.2144
$ws.Range("H27").Value = 6407.706
$ws.Range("I27").Value = 7497.95
$ws.Range("J27").Value = 4850.2144
$ws.Range("K27").Value = 7497.95
$ws.Range("M27").Value = -7390.95
$ws.Range("N27").Value = -5064.2144

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 31247.5
$ws.Range("I40").Value = 5000
$ws.Range("K40").Value = 5000
$ws.Range("M40").Value = -4851
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H70").Value = 85000
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 85000
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H81").Value = 19475.445
$ws.Range("I81").Value = 26296.5
$ws.Range("J81").Value = 5833.3335
$ws.Range("K81").Value = 52593
$ws.Range("L81").Value = 11666.667
$ws.Range("M81").Value = -51532
$ws.Range("N81").Value = -13788.667
$ws.Range("H84").Value = 19475.445
$ws.Range("I84").Value = 26296.5
$ws.Range("J84").Value = 5833.3335
$ws.Range("K84").Value = 262965
$ws.Range("L84").Value = 58333.335
$ws.Range("M84").Value = -257661
$ws.Range("N84").Value = -68941.33499999999
$ws.Range("H94").Value = 330000000
$ws.Range("J94").Value = 330000000
$ws.Range("L94").Value = 330000000
$ws.Range("N94").Value = -330001802
$ws.Range("H122").Value = 4571.391
$ws.Range("I122").Value = 3110.1333
$ws.Range("J122").Value = 7311.25
$ws.Range("K122").Value = 9330.3999
$ws.Range("L122").Value = 21933.75
$ws.Range("M122").Value = -6880.3999
$ws.Range("N122").Value = -26833.75
